$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Row 48/49 also swap the coin identity
# (RocketPoolETH and HuobiToken trade places) per the source diff.
$updates = @{
    "D2" = "42.414.66"
    "E2" = "  -0.64%  "
    "D3" = "2.289.74"
    "E3" = "  +0.28%  "
    "E4" = "  +0.05%  "
    "D5" = "301.63"
    "E5" = "  -1.10%  "
    "D6" = "95.70"
    "E6" = "  -0.44%  "
    "E7" = "  +0.11%  "
    "E8" = "  +0.12%  "
    "E9" = "  -1.90%  "
    "D10" = "34.44"
    "E10" = "  -2.34%  "
    "D11" = "19.03"
    "E11" = "  +3.66%  "
    "D12" = "0.0781"
    "E12" = "  -1.20%  "
    "E13" = "  +0.18%  "
    "E14" = "  +0.03%  "
    "D15" = "2.647.10"
    "E15" = "  +0.41%  "
    "D16" = "2.287.18"
    "E16" = "  +0.27%  "
    "E17" = "  +0.08%  "
    "D18" = "42.373.59"
    "E18" = "  -0.58%  "
    "E19" = "  -6.01%  "
    "D20" = "0.0₃0885"
    "E20" = "  -1.13%  "
    "E21" = "  -0.65%  "
    "D22" = "67.75"
    "E22" = "  +0.85%  "
    "E23" = "  +6.39%  "
    "D24" = "235.21"
    "E24" = "  -0.34%  "
    "E25" = "  -0.01%  "
    "D26" = "2.42"
    "E26" = "  -1.56%  "
    "D27" = "24.25"
    "E27" = "  -3.73%  "
    "E28" = "  +14.65%  "
    "D29" = "165.78"
    "E29" = "  -0.37%  "
    "E30" = "  -0.23%  "
    "E31" = "  -3.24%  "
    "E32" = "  +0.04%  "
    "E33" = "  +0.73%  "
    "D34" = "17.53"
    "E34" = "  -0.11%  "
    "D35" = "4.41"
    "E35" = "  -7.17%  "
    "D36" = "0.0698"
    "E36" = "  +0.84%  "
    "D38" = "0.0998"
    "E38" = "  -1.40%  "
    "E39" = "  -0.37%  "
    "E40" = "  -1.26%  "
    "E41" = "  -0.49%  "
    "D42" = "20.25"
    "E42" = "  +12.25%  "
    "D43" = "1.963.17"
    "E43" = "  -2.27%  "
    "D44" = "10.43"
    "E44" = "  +4.33%  "
    "D45" = "0.0278"
    "E45" = "  -0.19%  "
    "E46" = "  -0.79%  "
    "E47" = "  -0.98%  "
    "B48" = "HuobiToken"
    "C48" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D48" = "2.83"
    "E48" = "  -0.58%  "
    "B49" = "RocketPoolETH"
    "C49" = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
    "D49" = "2.517.03"
    "E49" = "  +0.54%  "
    "D50" = "53.10"
    "E50" = "  -1.08%  "
    "D51" = "71.06"
    "E51" = "  -0.41%  "
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force a text number format before assigning so Excel's automatic
    # type inference does not coerce numeric-looking strings (e.g. "1.00",
    # "42.414.66") into actual numbers or dates.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    # Restore the default "Normal" style so the cell keeps the same General
    # formatting / style index it had before this edit.
    $cell.Style = "Normal"
}
